$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 133
$ws.Range("F8").Value = 942
$ws.Range("F12").Value = 221
$ws.Range("F14").Value = 934
$ws.Range("F16").Value = 4053
$ws.Range("F17").Value = 1186
$ws.Range("F19").Value = 2645
$ws.Range("F21").Value = 1085
$ws.Range("F22").Value = 3620
$ws.Range("F23").Value = 774
$ws.Range("F26").Value = 2330
$ws.Range("F27").Value = 115
$ws.Range("F28").Value = 851
$ws.Range("F30").Value = 643
$ws.Range("F31").Value = 212
$ws.Range("F33").Value = 1365
$ws.Range("F34").Value = 1968
$ws.Range("F35").Value = 2
$ws.Range("F36").Value = 498
$ws.Range("F37").Value = 65
$ws.Range("F39").Value = 594
$ws.Range("F41").Value = 81
$ws.Range("F42").Value = 171
$ws.Range("F43").Value = 237
$ws.Range("F44").Value = 82

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F12").Value = 120

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 133
$ws.Range("F7").Value = 942
$ws.Range("F14").Value = 934
$ws.Range("F16").Value = 4053
$ws.Range("F17").Value = 1186
$ws.Range("F20").Value = 2645
$ws.Range("F21").Value = 1085
$ws.Range("F22").Value = 3620
$ws.Range("F23").Value = 774
$ws.Range("F27").Value = 2330
$ws.Range("F31").Value = 115
$ws.Range("F32").Value = 120
$ws.Range("F33").Value = 851
$ws.Range("F35").Value = 643
$ws.Range("F36").Value = 212
$ws.Range("F38").Value = 1365
$ws.Range("F39").Value = 1968
$ws.Range("F42").Value = 498
$ws.Range("F43").Value = 65
$ws.Range("F44").Value = 594
$ws.Range("F46").Value = 81
$ws.Range("F47").Value = 171
$ws.Range("F48").Value = 237
$ws.Range("F49").Value = 82
